$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N, shifting the existing
# "In Advance" (N) and "Over Due" (P) columns one to the right (-> O, Q).
$ws.Columns("N:N").Insert() | Out-Null

# Make "Repayment Schedule" the active sheet/tab and move the selection.
$ws.Activate() | Out-Null
$ws.Range("S5").Select() | Out-Null
